$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.823.77"
$ws.Range("E2").Value = "  -2.14%  "
$ws.Range("D3").Value = "2.306.86"
$ws.Range("E3").Value = "  -4.55%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'548.18"
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("D6").Value = "'131.57"
$ws.Range("E6").Value = "  -3.86%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -2.62%  "
$ws.Range("D9").Value = "2.305.12"
$ws.Range("E9").Value = "  -4.48%  "
$ws.Range("D10").Value = "'0.102"
$ws.Range("E10").Value = "  -2.92%  "
$ws.Range("D11").Value = "'5.55"
$ws.Range("E11").Value = "  -1.96%  "
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D13").Value = "'0.335"
$ws.Range("E13").Value = "  -5.22%  "
$ws.Range("D14").Value = "'24.03"
$ws.Range("E14").Value = "  -2.99%  "
$ws.Range("D15").Value = "2.722.57"
$ws.Range("E15").Value = "  -4.43%  "
$ws.Range("D16").Value = "58.803.26"
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("D17").Value = "'0.0000133"
$ws.Range("E17").Value = "  -3.25%  "
$ws.Range("D18").Value = "2.334.95"
$ws.Range("E18").Value = "  -3.25%  "
$ws.Range("D19").Value = "'10.71"
$ws.Range("E19").Value = "  -4.64%  "
$ws.Range("D20").Value = "'4.32"
$ws.Range("E20").Value = "  -4.61%  "
$ws.Range("D21").Value = "'314.90"
$ws.Range("E21").Value = "  -3.59%  "
$ws.Range("D22").Value = "'6.48"
$ws.Range("E22").Value = "  -4.77%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'63.34"
$ws.Range("E24").Value = "  -2.19%  "
$ws.Range("E25").Value = "  -4.96%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -5.99%  "
$ws.Range("D28").Value = "'1.31"
$ws.Range("E28").Value = "  -6.99%  "
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("D30").Value = "'169.53"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("D31").Value = "0.0₃0729"
$ws.Range("E31").Value = "  -5.67%  "
$ws.Range("D32").Value = "'1.11"
$ws.Range("E32").Value = "  +3.47%  "
$ws.Range("D33").Value = "'5.80"
$ws.Range("E33").Value = "  -5.05%  "
$ws.Range("D34").Value = "'0.384"
$ws.Range("E34").Value = "  -4.44%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'17.77"
$ws.Range("E36").Value = "  -3.80%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  -5.05%  "
$ws.Range("D39").Value = "'3.99"
$ws.Range("E39").Value = "  -5.71%  "
$ws.Range("D40").Value = "'38.14"
$ws.Range("E40").Value = "  -1.99%  "
$ws.Range("D41").Value = "'1.51"
$ws.Range("E41").Value = "  -5.19%  "
$ws.Range("D42").Value = "'301.11"
$ws.Range("E42").Value = "  -7.27%  "
$ws.Range("D43").Value = "'141.36"
$ws.Range("E43").Value = "  -3.45%  "
$ws.Range("D44").Value = "'3.44"
$ws.Range("E44").Value = "  -5.47%  "
$ws.Range("D45").Value = "'0.0953"
$ws.Range("E45").Value = "  -1.13%  "
$ws.Range("E46").Value = "  -2.57%  "
$ws.Range("D47").Value = "'0.557"
$ws.Range("E47").Value = "  -3.59%  "
$ws.Range("D48").Value = "'18.55"
$ws.Range("E48").Value = "  -6.28%  "
$ws.Range("E49").Value = "  -2.85%  "
$ws.Range("D50").Value = "'16.67"
$ws.Range("E50").Value = "  -4.28%  "
$ws.Range("D51").Value = "'11.01"
$ws.Range("E51").Value = "  -0.40%  "
